# Updated September List of GreyHat Presentations
# Reorders the 4 presentation rows (48-51) on the "GreyHat" section of the
# sheet so that the 3rd-Spring items (CS 4001, CS 4237...) come before the
# 4th-Fall items (CS 4210/4290..., CS 3210...), and updates the selected
# cell / scroll position to reflect where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Satisfied When?") is tied to the row's semester slot and is
# NOT part of the reshuffle - only the course details move between rows.
$cols = @("C", "D", "F", "G", "H", "I", "J")
$rows = @(48, 49, 50, 51)

# Snapshot the current ("before") values for rows 48-51, columns C,D,F..J
$snapshot = @{}
foreach ($r in $rows) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $snapshot[$addr] = $ws.Range($addr).Value()
    }
}

# New row N gets the contents of old row srcMap[N]
$srcMap = @{ 48 = 49; 49 = 51; 50 = 48; 51 = 50 }

foreach ($destRow in $rows) {
    $srcRow = $srcMap[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value = $snapshot[$srcAddr]
    }
}

# Update the sheet's selection to match where the editor left off
$ws.Range("C46").Select()
